$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.961.69'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '2.659.71'
$ws.Range('E3').Value = '  +4.00%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '514.70'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.27'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.568'
$ws.Range('E8').Value = '  +2.98%  '
$ws.Range('D9').Value = '2.691.22'
$ws.Range('E9').Value = '  +5.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.26'
$ws.Range('E10').Value = '  +1.02%  '
$ws.Range('E11').Value = '  +6.30%  '
$ws.Range('E12').Value = '  +1.99%  '
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('D14').Value = '3.132.84'
$ws.Range('E14').Value = '  +4.08%  '
$ws.Range('D15').Value = '58.951.62'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.06'
$ws.Range('E16').Value = '  +2.74%  '
$ws.Range('E17').Value = '  +2.71%  '
$ws.Range('D18').Value = '2.690.17'
$ws.Range('E18').Value = '  +4.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '348.64'
$ws.Range('E19').Value = '  +5.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.54'
$ws.Range('E20').Value = '  +1.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.43'
$ws.Range('E21').Value = '  +4.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.15'
$ws.Range('E22').Value = '  +3.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.00'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('E25').Value = '  +3.67%  '
$ws.Range('D26').Value = '2.793.33'
$ws.Range('E26').Value = '  +4.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0811'
$ws.Range('E29').Value = '  +4.77%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.28'
$ws.Range('E30').Value = '  +6.44%  '
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.45'
$ws.Range('E32').Value = '  +11.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.99'
$ws.Range('E33').Value = '  +2.53%  '
$ws.Range('E34').Value = '  +2.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '149.93'
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('E36').Value = '  +14.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.04'
$ws.Range('E37').Value = '  +3.14%  '
$ws.Range('E38').Value = '  +4.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.80'
$ws.Range('E39').Value = '  +3.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.846'
$ws.Range('E40').Value = '  +3.19%  '
$ws.Range('E41').Value = '  +6.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.40'
$ws.Range('E42').Value = '  +2.14%  '
$ws.Range('E43').Value = '  +2.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '279.00'
$ws.Range('E44').Value = '  -2.07%  '
$ws.Range('E45').Value = '  -0.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0983'
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.76'
$ws.Range('E47').Value = '  +6.22%  '
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0231'
$ws.Range('E49').Value = '  +2.63%  '
$ws.Range('D50').Value = '2.008.79'
$ws.Range('E50').Value = '  +5.34%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.72'
$ws.Range('E51').Value = '  +4.59%  '
